# Add a new code-snippet row (row 5) for "Convert an R file into an Rmd"
# and refresh the supporting formatting/view state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row content -------------------------------------------------
# Values are entered in the same order the shared-string table grew in
# the source edit: title, description, tags, url, src.
$ws.Range("A5").Value = 'Convert an `R` file into an Rmd'
$ws.Range("D5").Value = 'Use knitr::spin() to convert R file into Rmd'
$ws.Range("E5").Value = 'R; Convert to Rmd'
$ws.Range("C5").Value = 'https://github.com/sciencificity/convert-r-to-rmd'
$ws.Range("B5").Value = 'images/markus-spiske-hGb5WqRrWIg-unsplash.jpg'
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1

# --- Hyperlink + styling for the url cell -----------------------------
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/sciencificity/convert-r-to-rmd")
$ws.Range("C5").Style = "Hyperlink"

# --- Column width adjustments -----------------------------------------
$ws.Columns.Item(1).ColumnWidth = 27.6
$ws.Columns.Item(2).ColumnWidth = 33.6
$ws.Columns.Item(3).ColumnWidth = 57.6
$ws.Columns.Item(4).ColumnWidth = 22.3
$ws.Columns.Item(5).ColumnWidth = 28.6

# --- Selection / view state ---------------------------------------------
[void]$ws.Range("B5").Select()

# --- Workbook-level calculation option (iterative calculation) ---------
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.001
